# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
# Both sheets carry identical source data, and the commit regenerated the
# site's data export, bumping a handful of F-column counters.

$wb = $excel.ActiveWorkbook

$updates = @{
    3  = 771
    6  = 4550
    7  = 24
    8  = 367
    9  = 1314
    10 = 541
    11 = 53
    12 = 909
    14 = 507
}

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
